# sua chiet khau cua sale phu va update chien luoc chay tinh luong theo gio
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

# Ngay cong (working days) tang tu 19 len 21
$ws.Range("B2").Value = 21

# Phu cap = Ngay cong * 35,000 VND/ngay
$ws.Range("B3").Value = 735000

# Luong co ban tai LONG XUYEN = Ngay cong * (3,000,000 / 28)
$ws.Range("B12").Value = 2250000

# Tong luong tai LONG XUYEN = Luong co ban + Phu cap + Chiet khau sale chinh + Cong phu phau 1
$ws.Range("B29").Value = 3185000

# Tong luong (tat ca chi nhanh)
$ws.Range("B31").Value = 3185000
